# Add the "OrderStatusSearch++" test-data sheet (magento order status search)
# and wire it into the IC sheet's summary row, mirroring the author's
# "adding magento order search" commit.

$wb = $excel.ActiveWorkbook

$ic = $wb.Worksheets.Item("IC")

# New sheet goes right after "IC" and before "ProductSearch++".
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ic)
$newSheet.Name = "OrderStatusSearch++"

# IC sheet gets a new column (L) referencing the new suite by name.
$ic.Range("L2").Value = "OrderStatusSearch"

# Populate the new sheet's test data. Order matches the shared-string
# allocation sequence: D2, D1, C1, then the rest.
$newSheet.Range("D2").Value = "Canceled"
$newSheet.Range("D1").Value = "orderStatus"
$newSheet.Range("C1").Value = "productSearchId"
$newSheet.Range("A1").Value = "TCID"
$newSheet.Range("B1").Value = "occurence"
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = 1
$newSheet.Range("C2").Value = 3100000680

# Column widths (best-fit approximations) for the new sheet.
$newSheet.Columns.Item(2).ColumnWidth = 9.166666666666666
$newSheet.Columns.Item(3).ColumnWidth = 14.666666666666666
$newSheet.Columns.Item(4).ColumnWidth = 10.5

# IC sheet column L widens to fit the new header text.
$ic.Columns.Item(12).ColumnWidth = 16.833333333333336

# Restore selections: IC keeps a multi-row selection, new sheet gets its own.
$ic.Activate()
$ic.Range("A5:A26").EntireRow.Select() | Out-Null

$newSheet.Activate()
$newSheet.Range("E10").Select() | Out-Null
